$d = $word.ActiveDocument

$oldTitleText = "Play Deep Sea for Free - Review and Ratings"
$oldDescText  = "Explore Deep Sea, a beautifully designed underwater-themed online slot game. Read this review and play it for free."
$descSuffix   = ": " + $oldDescText

# ------------------------------------------------------------------
# Locate (by content, not a hard-coded index) the existing duplicate
# "Play Deep Sea for Free..." bold paragraph near the bottom of the
# document. We reuse its run structure (it already has the leading
# empty run + single bold run pattern used throughout this file) as
# the template for the new "Meta description" paragraph we're about
# to insert right under the Heading 1 title.
# ------------------------------------------------------------------
$sourceIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs($i).Range.Text.TrimEnd([char]13)
    if ($i -ne 1 -and $txt -eq $oldTitleText) {
        $sourceIndex = $i
        break
    }
}

$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

# The paragraph break we just inserted shifts every following
# paragraph (including our located source) down by one.
$metaSourcePara = $d.Paragraphs($sourceIndex + 1)

$newPara = $d.Paragraphs(2)
$newPara.Style = "Normal"
$newPara.Range.FormattedText = $metaSourcePara.Range.FormattedText

# newPara now reads "Play Deep Sea for Free - Review and Ratings" in a
# single bold run (plus the leading empty run). Overwrite just that
# run's text (bounded range, excludes the paragraph mark) with
# "Meta description".
$np = $d.Paragraphs(2)
$boldStart = $np.Range.Start
$boldRange = $d.Range($boldStart, $boldStart + $oldTitleText.Length)
$boldRange.Text = "Meta description"

# Append the (non-bold) description suffix as its own run right after
# the bold run, still inside the same paragraph.
$np2 = $d.Paragraphs(2)
$afterBold = $np2.Range.End - 1   # exclude the paragraph mark
$insPoint = $d.Range($afterBold, $afterBold)
$insPoint.InsertAfter($descSuffix)
$suffixRange = $d.Range($afterBold, $afterBold + $descSuffix.Length)
$suffixRange.Font.Bold = 0

# ------------------------------------------------------------------
# Remove the old duplicate "Play Deep Sea for Free..." paragraph near
# the bottom of the document (its content now lives at the top, right
# below the heading).
# ------------------------------------------------------------------
$total = $d.Paragraphs.Count
for ($i = $total; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    $txt = $p.Range.Text.TrimEnd([char]13)
    if ($i -ne 2 -and $txt -eq $oldTitleText) {
        $p.Range.Delete()
        break
    }
}

# ------------------------------------------------------------------
# Replace the final italic paragraph's text (the meta-description
# sentence) with the DALLE image prompt, keeping its italic run and
# its straight (non-curly) quote marks.
# ------------------------------------------------------------------
$q = [char]34
$dallePrompt = "Prompt for DALLE: Create a feature image for the game " + $q + "Deep Sea" + $q + " that depicts a happy Maya warrior with glasses. The image should be in a cartoon style, with bright and bold colors that reflect the underwater theme of the game. The warrior should be shown diving into the ocean, surrounded by sea creatures and marine beauty. The image should capture the excitement and adventure of the game, with a focus on the joy and energy of the Maya warrior as they explore the deep sea. Use bold lines and bright colors to make the image stand out and capture the attention of players. The image should be dynamic and eye-catching, drawing players in and encouraging them to dive into the world of " + $q + "Deep Sea." + $q

$descIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs($i).Range.Text.TrimEnd([char]13)
    if ($txt -eq $oldDescText) {
        $descIndex = $i
        break
    }
}
$descPara = $d.Paragraphs($descIndex)
$descStart = $descPara.Range.Start
$descRange = $d.Range($descStart, $descStart + $oldDescText.Length)
$descRange.Text = $dallePrompt

Write-Output "Done"
